$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-02 Monday" "2024-09-03 Tuesday"

Replace-Text "29÷7=" "23÷2="
Replace-Text "31÷7=" "24÷7="
Replace-Text "52÷2=" "27÷3="
Replace-Text "19÷7=" "67÷6="
Replace-Text "13÷4=" "19÷2="
Replace-Text "29÷3=" "75÷4="
Replace-Text "63÷8=" "48÷5="
Replace-Text "93÷8=" "10÷5="
Replace-Text "60÷8=" "64÷7="
Replace-Text "33÷5=" "61÷5="
Replace-Text "62÷8=" "47÷9="
Replace-Text "32÷6=" "94÷2="
Replace-Text "23÷3=" "39÷7="
Replace-Text "49÷3=" "51÷3="
Replace-Text "57÷9=" "20÷2="
Replace-Text "93÷9=" "99÷3="
Replace-Text "33÷7=" "54÷5="
Replace-Text "45÷4=" "68÷7="
Replace-Text "18÷7=" "29÷8="
Replace-Text "11÷9=" "78÷4="
Replace-Text "50÷7=" "95÷2="
Replace-Text "44÷3=" "93÷4="
Replace-Text "23÷8=" "10÷6="
Replace-Text "92÷6=" "77÷4="
Replace-Text "77÷2=" "55÷7="

Write-Output "Done applying replacements"
